# Add a new "update post" test-case row (row 8) to the Trainer sheet,
# re-using the login credentials already stored in row 1 (trainerict / Train123).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trainer")

# New test name for this case
$ws.Range("A8").Value = "Krishna Priya"

# Email column: copy the formatting used for the email in row 1, then set the value
$ws.Range("A1").Copy()
$ws.Range("B8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B8").Value = "trainerict@gmail.com"

# Password column: copy the formatting used for the password in row 1, then set the value
$ws.Range("B1").Copy()
$ws.Range("C8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C8").Value = "Train123"
